$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block (rows 14-26) ---

# Row 14
$ws.Range("A14").Value = 1000
$ws.Range("B14").Value = 1000000

# Labels in column A (rows 16-18, then 21-25) — filled in this order so the
# shared-string table indices line up with how Excel originally wrote them.
$ws.Range("A16").Value = "R"
$ws.Range("A17").Value = "Cp"
$ws.Range("A18").Value = "T_0"

$ws.Range("A21").Value = "a"
$ws.Range("A22").Value = "b"
$ws.Range("A23").Value = "T"
$ws.Range("A24").Value = "da"
$ws.Range("A25").Value = "db"

# Labels in column G (rows 14-16)
$ws.Range("G14").Value = "m1"
$ws.Range("G15").Value = "m2"
$ws.Range("G16").Value = "m3"

# Label in B20
$ws.Range("B20").Value = "T=24"

# Labels in column I (rows 14-16)
$ws.Range("I14").Value = "dm1"
$ws.Range("I15").Value = "dm2"
$ws.Range("I16").Value = "dm3"

# Final label, reuses the existing "dT" shared string
$ws.Range("A26").Value = "dT"

# Numeric / formula values

$ws.Range("H14").Formula = '=0.9242/$B$10'
$ws.Range("J14").Formula = '=H14*H11'

$ws.Range("H15").Formula = '=0.7537/$B$10'
$ws.Range("J15").Formula = '=H15*H12'

$ws.Range("B16").Value = 8.1300000000000008
$ws.Range("H16").Formula = '=0.6071/$B$10'
$ws.Range("J16").Formula = '=H16*H13'

$ws.Range("B17").Formula = '=(5/2+1)*$B$16'

$ws.Range("B18").Value = 273.16000000000003

$ws.Range("B21").Formula = '=$B$16*$B$17*($H$14-$H$15)/(1/($B$18+24)-1/($B$18+50))'
$ws.Range("B22").Formula = '=$B$17*($H$14*($B$18+24)-$H$15*($B$18+50))/(50-24)'
$ws.Range("B23").Formula = '=(2*$B$21)/($B$16*$B$22)'
$ws.Range("B24").Formula = '=$B$21/($H$14-$H$15)*SQRT($J$14^2+$J$15^2)'
$ws.Range("B25").Formula = '=$B$17*SQRT($J$14*($B$18+24)^2-$J$15*($B$18+50)^2)/(50-24)'

# --- Column widths (best-fit) for the new H/J columns ---
# Set BEFORE repositioning the chart, since the chart's pixel position is
# re-derived from the (now different) column widths.
$ws.Columns.Item(8).ColumnWidth = 11.17
$ws.Columns.Item(10).ColumnWidth = 11.17

# --- Move / resize the chart to make room for the new data ---
$co = $ws.ChartObjects(1)
$co.Left = 713.77734375
$co.Top = 147.37496062992125
$co.Width = 442.0625
$co.Height = 216

# --- View state: scroll position + selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D13").Select()
